$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: EAP moves away from F3
$ws.Range("F3").Value = "-"

# Row 6: EAP moves from C6 to E6
$ws.Range("C6").Value = "-"
$ws.Range("E6").Value = "EAP"
